# Auto-generated edit script: Add data for 2022-12-06
# Updates column I (2022 totals) across 38 worksheets to reflect one additional day of data.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Citywide Totals")
$ws.Cells.Item(2, 9).Value = 6829
$ws.Cells.Item(3, 9).Value = 7099
$ws.Cells.Item(4, 9).Value = 1629
$ws.Cells.Item(5, 9).Value = 661
$ws.Cells.Item(6, 9).Value = 8269
$ws.Cells.Item(7, 9).Value = 24487

$ws = $wb.Worksheets.Item("Grand Crossing")
$ws.Cells.Item(2, 9).Value = 224
$ws.Cells.Item(7, 9).Value = 759

$ws = $wb.Worksheets.Item("Woodlawn")
$ws.Cells.Item(6, 9).Value = 110
$ws.Cells.Item(7, 9).Value = 431

$ws = $wb.Worksheets.Item("North Lawndale")
$ws.Cells.Item(2, 9).Value = 220
$ws.Cells.Item(3, 9).Value = 347
$ws.Cells.Item(5, 9).Value = 29
$ws.Cells.Item(6, 9).Value = 279
$ws.Cells.Item(7, 9).Value = 931

$ws = $wb.Worksheets.Item("New City")
$ws.Cells.Item(2, 9).Value = 184
$ws.Cells.Item(6, 9).Value = 175
$ws.Cells.Item(7, 9).Value = 570

$ws = $wb.Worksheets.Item("By Neighborhood")
$ws.Cells.Item(2, 9).Value = 192
$ws.Cells.Item(3, 9).Value = 17
$ws.Cells.Item(7, 9).Value = 768
$ws.Cells.Item(8, 9).Value = 1464
$ws.Cells.Item(9, 9).Value = 127
$ws.Cells.Item(11, 9).Value = 374
$ws.Cells.Item(18, 9).Value = 194
$ws.Cells.Item(19, 9).Value = 687
$ws.Cells.Item(25, 9).Value = 128
$ws.Cells.Item(27, 9).Value = 213
$ws.Cells.Item(29, 9).Value = 1477
$ws.Cells.Item(33, 9).Value = 1086
$ws.Cells.Item(34, 9).Value = 112
$ws.Cells.Item(36, 9).Value = 333
$ws.Cells.Item(37, 9).Value = 759
$ws.Cells.Item(40, 9).Value = 43
$ws.Cells.Item(41, 9).Value = 107
$ws.Cells.Item(42, 9).Value = 904
$ws.Cells.Item(44, 9).Value = 186
$ws.Cells.Item(48, 9).Value = 309
$ws.Cells.Item(55, 9).Value = 282
$ws.Cells.Item(60, 9).Value = 140
$ws.Cells.Item(61, 9).Value = 26
$ws.Cells.Item(63, 9).Value = 74
$ws.Cells.Item(64, 9).Value = 195
$ws.Cells.Item(65, 9).Value = 570
$ws.Cells.Item(67, 9).Value = 931
$ws.Cells.Item(73, 9).Value = 219
$ws.Cells.Item(77, 9).Value = 146
$ws.Cells.Item(78, 9).Value = 330
$ws.Cells.Item(79, 9).Value = 703
$ws.Cells.Item(83, 9).Value = 528
$ws.Cells.Item(85, 9).Value = 1096
$ws.Cells.Item(86, 9).Value = 157
$ws.Cells.Item(90, 9).Value = 315
$ws.Cells.Item(91, 9).Value = 258
$ws.Cells.Item(99, 9).Value = 431
$ws.Cells.Item(101, 9).Value = 24487

$ws = $wb.Worksheets.Item("South Chicago")
$ws.Cells.Item(2, 9).Value = 176
$ws.Cells.Item(6, 9).Value = 118
$ws.Cells.Item(7, 9).Value = 528

$ws = $wb.Worksheets.Item("Garfield Park")
$ws.Cells.Item(2, 9).Value = 245
$ws.Cells.Item(3, 9).Value = 400
$ws.Cells.Item(6, 9).Value = 349
$ws.Cells.Item(7, 9).Value = 1086

$ws = $wb.Worksheets.Item("Englewood")
$ws.Cells.Item(2, 9).Value = 436
$ws.Cells.Item(3, 9).Value = 507
$ws.Cells.Item(4, 9).Value = 78
$ws.Cells.Item(6, 9).Value = 408
$ws.Cells.Item(7, 9).Value = 1477

$ws = $wb.Worksheets.Item("Chatham")
$ws.Cells.Item(2, 9).Value = 222
$ws.Cells.Item(3, 9).Value = 202
$ws.Cells.Item(6, 9).Value = 219
$ws.Cells.Item(7, 9).Value = 687

$ws = $wb.Worksheets.Item("Irving Park")
$ws.Cells.Item(2, 9).Value = 61
$ws.Cells.Item(7, 9).Value = 186

$ws = $wb.Worksheets.Item("Lake View")
$ws.Cells.Item(6, 9).Value = 160
$ws.Cells.Item(7, 9).Value = 309

$ws = $wb.Worksheets.Item("South Shore")
$ws.Cells.Item(2, 9).Value = 314
$ws.Cells.Item(3, 9).Value = 413
$ws.Cells.Item(7, 9).Value = 1096

$ws = $wb.Worksheets.Item("Hermosa")
$ws.Cells.Item(6, 9).Value = 30
$ws.Cells.Item(7, 9).Value = 107

$ws = $wb.Worksheets.Item("Humboldt Park")
$ws.Cells.Item(2, 9).Value = 208
$ws.Cells.Item(3, 9).Value = 270
$ws.Cells.Item(6, 9).Value = 343
$ws.Cells.Item(7, 9).Value = 904

$ws = $wb.Worksheets.Item("Rogers Park")
$ws.Cells.Item(2, 9).Value = 74
$ws.Cells.Item(4, 9).Value = 45
$ws.Cells.Item(6, 9).Value = 119
$ws.Cells.Item(7, 9).Value = 330

$ws = $wb.Worksheets.Item("Lower West Side")
$ws.Cells.Item(2, 9).Value = 86
$ws.Cells.Item(7, 9).Value = 282

$ws = $wb.Worksheets.Item("Washington Park")
$ws.Cells.Item(2, 9).Value = 81
$ws.Cells.Item(7, 9).Value = 258

$ws = $wb.Worksheets.Item("Roseland")
$ws.Cells.Item(2, 9).Value = 207
$ws.Cells.Item(3, 9).Value = 228
$ws.Cells.Item(6, 9).Value = 201
$ws.Cells.Item(7, 9).Value = 703

$ws = $wb.Worksheets.Item("Near South Side")
$ws.Cells.Item(6, 9).Value = 65
$ws.Cells.Item(7, 9).Value = 195

$ws = $wb.Worksheets.Item("Calumet Heights")
$ws.Cells.Item(4, 9).Value = 7
$ws.Cells.Item(7, 9).Value = 194

$ws = $wb.Worksheets.Item("Grand Boulevard")
$ws.Cells.Item(6, 9).Value = 104
$ws.Cells.Item(7, 9).Value = 333

$ws = $wb.Worksheets.Item("Garfield Ridge")
$ws.Cells.Item(2, 9).Value = 45
$ws.Cells.Item(3, 9).Value = 32
$ws.Cells.Item(7, 9).Value = 112

$ws = $wb.Worksheets.Item("East Side")
$ws.Cells.Item(2, 9).Value = 48
$ws.Cells.Item(7, 9).Value = 128

$ws = $wb.Worksheets.Item("Belmont Cragin")
$ws.Cells.Item(2, 9).Value = 148
$ws.Cells.Item(6, 9).Value = 100
$ws.Cells.Item(7, 9).Value = 374

$ws = $wb.Worksheets.Item("Avalon Park")
$ws.Cells.Item(6, 9).Value = 39
$ws.Cells.Item(7, 9).Value = 127

$ws = $wb.Worksheets.Item("Portage Park")
$ws.Cells.Item(3, 9).Value = 67
$ws.Cells.Item(7, 9).Value = 219

$ws = $wb.Worksheets.Item("Albany Park")
$ws.Cells.Item(2, 9).Value = 72
$ws.Cells.Item(3, 9).Value = 61
$ws.Cells.Item(7, 9).Value = 192

$ws = $wb.Worksheets.Item("Austin")
$ws.Cells.Item(2, 9).Value = 435
$ws.Cells.Item(3, 9).Value = 422
$ws.Cells.Item(5, 9).Value = 47
$ws.Cells.Item(6, 9).Value = 469
$ws.Cells.Item(7, 9).Value = 1464

$ws = $wb.Worksheets.Item("Edgewater")
$ws.Cells.Item(6, 9).Value = 83
$ws.Cells.Item(7, 9).Value = 213

$ws = $wb.Worksheets.Item("Streeterville")
$ws.Cells.Item(2, 9).Value = 29
$ws.Cells.Item(4, 9).Value = 75
$ws.Cells.Item(6, 9).Value = 36
$ws.Cells.Item(7, 9).Value = 157

$ws = $wb.Worksheets.Item("Washington Heights")
$ws.Cells.Item(3, 9).Value = 80
$ws.Cells.Item(7, 9).Value = 315

$ws = $wb.Worksheets.Item("Morgan Park")
$ws.Cells.Item(2, 9).Value = 49
$ws.Cells.Item(7, 9).Value = 140

$ws = $wb.Worksheets.Item("Riverdale")
$ws.Cells.Item(6, 9).Value = 37
$ws.Cells.Item(7, 9).Value = 146

$ws = $wb.Worksheets.Item("Andersonville")
$ws.Cells.Item(4, 9).Value = 2
$ws.Cells.Item(6, 9).Value = 17

$ws = $wb.Worksheets.Item("Hegewisch")
$ws.Cells.Item(2, 9).Value = 16
$ws.Cells.Item(7, 9).Value = 43

$ws = $wb.Worksheets.Item("Auburn Gresham")
$ws.Cells.Item(6, 9).Value = 208
$ws.Cells.Item(7, 9).Value = 768

$ws = $wb.Worksheets.Item("Mount Greenwood")
$ws.Cells.Item(6, 9).Value = 11
$ws.Cells.Item(7, 9).Value = 26
